# "cap nhat canh buy" - add a new highlighted note cell (J12) to the plan.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new note text in J12 (this also grows the new shared string entry).
$cell = $ws.Range("J12")
$cell.Value = "Lĩnh ngộ tinh túy thị trường"

# Give it the same Bold/Calibri look already used for the "RSI LÀ GÌ" heading
# (A14), then bump the size up to 20pt so it stands out - this reuses the
# existing font/style catalog entries as much as possible and only adds the
# one new Bold/20pt/Calibri font+style pairing that the workbook needs.
$heading = $ws.Range("A14")
$heading.Copy() | Out-Null
$cell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$cell.Font.Size = 20
$excel.CutCopyMode = $false | Out-Null

# The taller note text means row 12 needs to be taller too.
$ws.Rows.Item(12).RowHeight = 26.25

# Leave the view scrolled back to the top, with the new cell's neighbour
# selected, instead of the far-away cell that was selected before.
$ws.Range("I12").Select() | Out-Null
